# Delete the SMTP/IMAP mail-server config rows from the "Config" sheet
# (rows 2-6: BotMailProtocol, Bot_SMTP_MailServer, Bot_SMTP_MailPort,
#  Bot_IMAP_MailServer, Bot_IMAP_MailPort), shifting the remaining rows
# (SystemExceptionEmailSubject, BussinessExceptionEmailSubject,
#  SystemExceptionEmailBody, BusinessExceptionEmailBody) up to rows 2-5.
# Also move the active/selected tab from "Constants" back to "Config".

$wb = $excel.ActiveWorkbook

$configSheet = $wb.Worksheets.Item("Config")
$configSheet.Rows("2:6").Delete()

# Make "Config" the selected/active sheet again (it was "Constants" before),
# and put the selection on the last data cell (B5), matching the new extent.
$configSheet.Select()
$configSheet.Range("B5").Select()
